$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "76.307.60"
$ws.Range("E2").Value = "  +0.44%  "
# Row 3
$ws.Range("D3").Value = "2.986.20"
$ws.Range("E3").Value = "  +2.15%  "
# Row 4
$ws.Range("E4").Value = "  +0.03%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "199.41"
$ws.Range("E5").Value = "  +0.52%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "627.52"
$ws.Range("E6").Value = "  +4.17%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.01%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.548"
$ws.Range("E8").Value = "  -1.14%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.199"
$ws.Range("E9").Value = "  +2.90%  "
# Row 10
$ws.Range("D10").Value = "2.986.24"
$ws.Range("E10").Value = "  +2.18%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.436"
$ws.Range("E11").Value = "  +6.01%  "
# Row 12
$ws.Range("E12").Value = "  -0.07%  "
# Row 13
$ws.Range("D13").Value = "3.522.92"
$ws.Range("E13").Value = "  +2.28%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.92"
$ws.Range("E14").Value = "  -0.20%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.78"
$ws.Range("E15").Value = "  +4.31%  "
# Row 16
$ws.Range("D16").Value = "76.269.26"
$ws.Range("E16").Value = "  +0.55%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000186"
$ws.Range("E17").Value = "  -1.74%  "
# Row 18
$ws.Range("D18").Value = "2.990.27"
$ws.Range("E18").Value = "  +2.98%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.44"
$ws.Range("E19").Value = "  +5.96%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.77"
$ws.Range("E20").Value = "  -1.68%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "372.91"
$ws.Range("E21").Value = "  -1.62%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.26"
$ws.Range("E22").Value = "  -2.84%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.30"
$ws.Range("E23").Value = "  +3.24%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.78"
$ws.Range("E24").Value = "  +1.18%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.09%  "
# Row 26
$ws.Range("D26").Value = "3.048.58"
$ws.Range("E26").Value = "  -0.52%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.27"
$ws.Range("E27").Value = "  +0.02%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.64"
$ws.Range("E28").Value = "  -0.67%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000107"
$ws.Range("E29").Value = "  -4.23%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("E30").Value = "  -1.04%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.31"
$ws.Range("E31").Value = "  +6.86%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.39"
$ws.Range("E32").Value = "  -2.03%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "503.07"
$ws.Range("E33").Value = "  -1.49%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.93"
$ws.Range("E34").Value = "  +5.21%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.02%  "
# Row 36
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "20.25"
$ws.Range("E36").Value = "  +0.43%  "
# Row 37
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "163.76"
$ws.Range("E37").Value = "  +0.47%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.386"
$ws.Range("E38").Value = "  +11.75%  "
# Row 39
$ws.Range("E39").Value = "  +1.40%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.105"
$ws.Range("E40").Value = "  +14.12%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.111"
$ws.Range("E41").Value = "  -3.27%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "182.55"
$ws.Range("E42").Value = "  +0.87%  "
# Row 43
$ws.Range("E43").Value = "  -0.02%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.93"
$ws.Range("E44").Value = "  -2.56%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.63"
$ws.Range("E45").Value = "  -3.11%  "
# Row 46
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.24"
$ws.Range("E46").Value = "  +0.11%  "
# Row 47
$ws.Range("B47").Value = "ImmutableX"
$ws.Range("C47").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.20"
$ws.Range("E47").Value = "  -2.96%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.590"
$ws.Range("E48").Value = "  +1.17%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.32"
$ws.Range("E49").Value = "  -4.49%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.83"
$ws.Range("E50").Value = "  +2.36%  "
# Row 51
$ws.Range("E51").Value = "  +1.25%  "
